$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "JSONMessageType" labels to "JSONType" in column J (rows 2, 6, 7)
$ws.Range("J2").Value = "JSONType"
$ws.Range("J6").Value = "JSONType"
$ws.Range("J7").Value = "JSONType"

# Keep selection on J7 as last active cell
$ws.Range("J7").Select()
